# Auto-generated edit script applying the cryptos.xlsx price-table refresh
# (GitHub Actions automated update, Thu Dec 28 18:56:12 UTC 2023).
#
# Column D ("Price") holds values that look numeric ("330.57", "42.417.54")
# but must stay plain text, matching the source inlineStr cells. Excel's
# Range.Value setter auto-converts numeric-looking strings to real numbers
# (losing formatting / multi-dot "thousands" groups), so for column D we:
#   1) force the cell to Text format ("@") before assigning,
#   2) assign the literal string,
#   3) reset the cell style back to "Normal" so no stray number format
#      sticks around on the saved cell (keeps the style identical to before).
# Column E (percentages, e.g. "  -1.72%  ") and columns B/C (coin name /
# link, including the Uniswap <-> WrappedBTC row swap) are plain text
# already, so a direct Value assignment is enough.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.582.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.357.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.638"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("E12").Value = "  -5.06%  "
$ws.Range("E13").Value = "  -4.36%  "
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.710.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.352.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.458.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.47%  "
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("E21").Value = "  +8.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.99%  "
$ws.Range("E24").Value = "  -10.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0901"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  -8.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.84%  "
$ws.Range("E38").Value = "  -5.02%  "
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "117.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +30.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.81%  "
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.566.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.34%  "
